$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reporting year (2022) was appended as column M, mirroring column L
# (2021) in both formatting (borders, number styles) and -- for the
# still-identical placeholder figures -- the same values as 2021.

# Copy column L's formatting (row 3, the bottom-border spacer row, through
# row 11, the footer row) into the new column M so the new column's borders
# and styles line up with the rest of the table.
$ws.Range("L3:L11").Copy()
$ws.Range("M3:M11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new column's own values: the header year reads 2022, and each
# data row repeats the same figure already shown for 2021.
$ws.Range("M4").Value = 2022
$ws.Range("M5").Value = 0.86
$ws.Range("M6").Value = 1.07
$ws.Range("M7").Value = 25.27
$ws.Range("M8").Value = 14
$ws.Range("M9").Value = 0.12
$ws.Range("M10").Value = 21.74
$ws.Range("M11").Value = 9.4600000000000009

# The active selection recorded for this edit moved to N6.
$ws.Range("N6").Select()
